$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 592.75
$ws.Range("J17").Value = 592.75
$ws.Range("L17").Value = 1778.25
$ws.Range("N17").Value = -2114.25
$ws.Range("H40").Value = 1890.0667
$ws.Range("I40").Value = 1945.381
$ws.Range("J40").Value = 1761
$ws.Range("K40").Value = 1945.381
$ws.Range("L40").Value = 1761
$ws.Range("M40").Value = -1770.381
$ws.Range("N40").Value = -2111
$ws.Range("H62").Value = 157145000
$ws.Range("I62").Value = 183335330
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 183335330
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -183334706
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 157145000
$ws.Range("I65").Value = 183335330
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 916676650
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -916673530
$ws.Range("N65").Value = -21240
$ws.Range("H96").Value = 8337855
$ws.Range("J96").Value = 20000454
$ws.Range("L96").Value = 60001362
$ws.Range("N96").Value = -60004108
$ws.Range("H97").Value = 6127.5
$ws.Range("J97").Value = 6127.5
$ws.Range("L97").Value = 18382.5
$ws.Range("N97").Value = -19374.5
$ws.Range("H100").Value = 3535.9546
$ws.Range("I100").Value = 2732.0833
$ws.Range("J100").Value = 4500.6
$ws.Range("K100").Value = 2732.0833
$ws.Range("L100").Value = 4500.6
$ws.Range("M100").Value = -2191.0833
$ws.Range("N100").Value = -5582.6
$ws.Range("H101").Value = 3196
$ws.Range("I101").Value = 4974.75
$ws.Range("J101").Value = 1773
$ws.Range("K101").Value = 14924.25
$ws.Range("L101").Value = 5319
$ws.Range("M101").Value = -13302.25
$ws.Range("N101").Value = -8563
$ws.Range("H137").Value = 4469403
$ws.Range("I137").Value = 9617384
$ws.Range("J137").Value = 7819.433
$ws.Range("K137").Value = 28852152
$ws.Range("L137").Value = 23458.299
$ws.Range("M137").Value = -28849602
$ws.Range("N137").Value = -28558.299

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 3385.5715
$ws.Range("I31").Value = 3385.5715
$ws.Range("K31").Value = 3385.5715
$ws.Range("M31").Value = -3091.5715
$ws.Range("H45").Value = 3830.75
$ws.Range("I45").Value = 3511.5
$ws.Range("K45").Value = 3511.5
$ws.Range("M45").Value = -3134.5
$ws.Range("H132").Value = 7824.967
$ws.Range("I132").Value = 4250.0625
$ws.Range("K132").Value = 12750.1875
$ws.Range("M132").Value = -10220.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 32364
$ws.Range("J21").Value = 32364
$ws.Range("L21").Value = 32364
$ws.Range("N21").Value = -32836
$ws.Range("H86").Value = 69247.97
$ws.Range("I86").Value = 2663.5715
$ws.Range("J86").Value = 224611.56
$ws.Range("K86").Value = 2663.5715
$ws.Range("L86").Value = 224611.56
$ws.Range("M86").Value = -1540.5715
$ws.Range("N86").Value = -226857.56
$ws.Range("H89").Value = 69247.97
$ws.Range("I89").Value = 2663.5715
$ws.Range("J89").Value = 224611.56
$ws.Range("K89").Value = 13317.8575
$ws.Range("L89").Value = 1123057.8
$ws.Range("M89").Value = -7701.8575
$ws.Range("N89").Value = -1134289.8
$ws.Range("H94").Value = 1115.76
$ws.Range("I94").Value = 1189.1578
$ws.Range("J94").Value = 883.3333
$ws.Range("K94").Value = 1189.1578
$ws.Range("L94").Value = 883.3333
$ws.Range("M94").Value = -738.1578
$ws.Range("N94").Value = -1785.3333
$ws.Range("H102").Value = 14368
$ws.Range("I102").Value = 10192.8
$ws.Range("K102").Value = 10192.8
$ws.Range("M102").Value = -6947.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 668.4643
$ws.Range("I7").Value = 329.53333
$ws.Range("J7").Value = 1059.5385
$ws.Range("K7").Value = 329.53333
$ws.Range("L7").Value = 1059.5385
$ws.Range("M7").Value = -216.53333
$ws.Range("N7").Value = -1285.5385
$ws.Range("H28").Value = 27940.166
$ws.Range("J28").Value = 27940.166
$ws.Range("L28").Value = 27940.166
$ws.Range("N28").Value = -28430.166
$ws.Range("H31").Value = 26318986
$ws.Range("J31").Value = 4549.4287
$ws.Range("L31").Value = 4549.4287
$ws.Range("N31").Value = -5139.4287
$ws.Range("H34").Value = 26318986
$ws.Range("J34").Value = 4549.4287
$ws.Range("L34").Value = 4549.4287
$ws.Range("N34").Value = -4953.4287
$ws.Range("H62").Value = 19845
$ws.Range("I62").Value = 17497.75
$ws.Range("K62").Value = 17497.75
$ws.Range("M62").Value = -16873.75
$ws.Range("H65").Value = 19845
$ws.Range("I65").Value = 17497.75
$ws.Range("K65").Value = 87488.75
$ws.Range("M65").Value = -84368.75
$ws.Range("H134").Value = 4857.2
$ws.Range("I134").Value = 3949.0977
$ws.Range("K134").Value = 11847.2931
$ws.Range("M134").Value = -9312.293099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 125.190475
$ws.Range("I33").Value = 73.625
$ws.Range("K33").Value = 441.75
$ws.Range("M33").Value = -158.75
$ws.Range("H40").Value = 842.8
$ws.Range("I40").Value = 106.85714
$ws.Range("K40").Value = 427.42856
$ws.Range("M40").Value = -358.42856
$ws.Range("H44").Value = 2054.2
$ws.Range("J44").Value = 2890.5
$ws.Range("L44").Value = 8671.5
$ws.Range("N44").Value = -9467.5
$ws.Range("H68").Value = 53226.05
$ws.Range("J68").Value = 3510.8438
$ws.Range("L68").Value = 10532.5314
$ws.Range("N68").Value = -12154.5314
$ws.Range("H71").Value = 53226.05
$ws.Range("J71").Value = 3510.8438
$ws.Range("L71").Value = 31597.5942
$ws.Range("N71").Value = -39709.5942
$ws.Range("H86").Value = 304.72726
$ws.Range("I86").Value = 285.3
$ws.Range("J86").Value = 499
$ws.Range("K86").Value = 855.9000000000001
$ws.Range("L86").Value = 1497
$ws.Range("M86").Value = 330.0999999999999
$ws.Range("N86").Value = -3869
$ws.Range("H89").Value = 304.72726
$ws.Range("I89").Value = 285.3
$ws.Range("J89").Value = 499
$ws.Range("K89").Value = 2567.7
$ws.Range("L89").Value = 4491
$ws.Range("M89").Value = 3360.3
$ws.Range("N89").Value = -16347
$ws.Range("H121").Value = 16668284
$ws.Range("I121").Value = 364.5
$ws.Range("J121").Value = 18520276
$ws.Range("K121").Value = 1093.5
$ws.Range("L121").Value = 55560828
$ws.Range("M121").Value = 216.5
$ws.Range("N121").Value = -55563448
$ws.Range("H137").Value = 2836.125
$ws.Range("I137").Value = 1457.8
$ws.Range("J137").Value = 5133.3335
$ws.Range("K137").Value = 4373.4
$ws.Range("L137").Value = 15400.0005
$ws.Range("M137").Value = 726.6000000000004
$ws.Range("N137").Value = -25600.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4323.5347
$ws.Range("I122").Value = 4331.2383
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 12993.7149
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -10543.7149
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 4550.1724
$ws.Range("I132").Value = 2167.8
$ws.Range("K132").Value = 6503.400000000001
$ws.Range("M132").Value = -3973.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 16132338
$ws.Range("I122").Value = 17860392
$ws.Range("K122").Value = 53581176
$ws.Range("M122").Value = -53578726
$ws.Range("H132").Value = 4201.091
$ws.Range("I132").Value = 3477.1035
$ws.Range("K132").Value = 10431.3105
$ws.Range("M132").Value = -7901.3105

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 8336357
$ws.Range("I126").Value = 10002708
$ws.Range("K126").Value = 30008124
$ws.Range("M126").Value = -30005654
